$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:N) to (B:O)
$ws.Columns.Item(1).Insert(-4161)

# Fill new column A with fold labels for rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "fold " + ($r - 1)
}

# Apply the same formatting used for the header row (bold, centered, bordered)
$rng = $ws.Range("A2:A11")
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

Write-Output "done"
